$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of professor data (NOM in column B, PRENOM in column C)
$ws.Range("B2").Value = "El Haddad"
$ws.Range("C2").Value = "Mohamed"

$ws.Range("B3").Value = "Hassan"
$ws.Range("C3").Value = "Badir"

$ws.Range("B4").Value = "El Alami"
$ws.Range("C4").Value = "Hassoun"

$ws.Range("B5").Value = "Amechnoue"
$ws.Range("C5").Value = "Khalid"

$ws.Range("B6").Value = "Fissoune"
$ws.Range("C6").Value = "Rachida"

$ws.Range("B7").Value = "Azzouzi"
$ws.Range("C7").Value = "Rahali"

# Header for new column E (added last so it lands at the end of sharedStrings)
$ws.Range("E1").Value = "Matières enseignés"

# Column widths (closest achievable values to the target stored widths
# 33.7109375 / 31.85546875, given this engine's column-width rounding model)
$ws.Columns.Item(4).ColumnWidth = 32.833333333333336
$ws.Columns.Item(5).ColumnWidth = 31.0

# Selection matches diff (activeCell D9)
$ws.Range("D9").Select()
